# Handle read and write file excel
# Fix up the sample data in both sheets (rows 3-12 of column A should count
# on up from 3 instead of repeating the 1/2 pattern) and correct which
# sheet/cell is left selected when the file is saved.

$wb = $excel.ActiveWorkbook

$basic = $wb.Worksheets.Item("Basic")
$advance = $wb.Worksheets.Item("Advance")

# --- Re-number column A (rows 3-12) on the "Basic" sheet ---
$basic.Range("A3").Value = 3
$basic.Range("A4").Value = 4
$basic.Range("A5").Value = 5
$basic.Range("A6").Value = 6
$basic.Range("A7").Value = 7
$basic.Range("A8").Value = 8
$basic.Range("A9").Value = 9
$basic.Range("A10").Value = 10
$basic.Range("A11").Value = 11
$basic.Range("A12").Value = 12

# --- Re-number column A (rows 3-12) on the "Advance" sheet ---
$advance.Range("A3").Value = 3
$advance.Range("A4").Value = 4
$advance.Range("A5").Value = 5
$advance.Range("A6").Value = 6
$advance.Range("A7").Value = 7
$advance.Range("A8").Value = 8
$advance.Range("A9").Value = 9
$advance.Range("A10").Value = 10
$advance.Range("A11").Value = 11
$advance.Range("A12").Value = 12

# "Advance" is the active sheet coming in; leave its cursor on C18 before
# switching away from it.
$advance.Range("C18").Select()

# Make "Basic" the active/selected sheet and park the cursor on E16, which
# becomes the saved selection state for the workbook.
$basic.Activate()
$basic.Range("E16").Select()
